# Update header labels in the plate import template.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "compound"
$ws.Range("C3").Value = "compound name"
$ws.Range("A31").Value = "readout"
$ws.Range("A41").Value = "read_norm"
$ws.Range("A11").Value = "conc"

# Move the active selection, matching the saved cursor position.
$ws.Range("A12").Select()

$wb.Save()
